$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.311.38"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "3.100.60"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'579.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "'171.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.63%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.096.54"
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("D9").Value = "'0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("D10").Value = "'6.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("E11").Value = "  -2.40%  "

$ws.Range("D12").Value = "'0.477"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("E13").Value = "  -1.96%  "

$ws.Range("D14").Value = "'36.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.44%  "

$ws.Range("D15").Value = "'0.122"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.81%  "

$ws.Range("D16").Value = "3.615.61"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").Value = "67.188.26"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "'7.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("D19").Value = "'16.67"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "3.094.68"
$ws.Range("E20").Value = "  -1.37%  "

$ws.Range("D21").Value = "'488.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").Value = "'7.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").Value = "'0.696"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.87%  "

$ws.Range("D24").Value = "'83.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.90%  "

$ws.Range("D25").Value = "'13.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.95%  "

$ws.Range("D26").Value = "'2.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "

$ws.Range("D27").Value = "'10.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.80%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").Value = "'7.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.64%  "

$ws.Range("D30").Value = "'2.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.63%  "

$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("D32").Value = "'28.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.42%  "

$ws.Range("E33").Value = "  -1.53%  "

$ws.Range("D34").Value = "0.0₃0943"
$ws.Range("E34").Value = "  -5.66%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").Value = "'5.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("D38").Value = "'46.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("E39").Value = "  -4.56%  "

$ws.Range("D40").Value = "'0.124"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.70%  "

$ws.Range("E41").Value = "  -2.42%  "

$ws.Range("D42").Value = "'8.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.80%  "

$ws.Range("D43").Value = "2.792.75"
$ws.Range("E43").Value = "  -1.95%  "

$ws.Range("D44").Value = "'380.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("E45").Value = "  -7.51%  "

$ws.Range("D46").Value = "'0.0350"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("D47").Value = "'135.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D49").Value = "'24.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("E51").Value = "  -1.69%  "
